# Update the panel-query timestamps on the "data" sheet (F2:F8) to reflect
# the new fetch run.
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$data.Range("F2").Value = "2021-10-05 14:19:42.640651"
$data.Range("F3").Value = "2021-10-05 14:19:42.640658"
$data.Range("F4").Value = "2021-10-05 14:19:42.640661"
$data.Range("F5").Value = "2021-10-05 14:19:42.640664"
$data.Range("F6").Value = "2021-10-05 14:19:42.640667"
$data.Range("F7").Value = "2021-10-05 14:19:42.640670"
$data.Range("F8").Value = "2021-10-05 14:19:42.640672"

# Add a new "metadata" tab after the existing "data" tab.
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $data)
$meta.Name = "metadata"

# Give the metadata header row (B1:G1) and the index cell (A2) the same
# bold/bordered/centered style used for the header row & index column on
# the "data" sheet, by copying the format from an already-styled cell.
$headerCols = @("B", "C", "D", "E", "F", "G")
foreach ($col in $headerCols) {
    $data.Range("B1").Copy()
    $meta.Range($col + "1").PasteSpecial(-4122)
}
$data.Range("B1").Copy()
$meta.Range("A2").PasteSpecial(-4122)

# Header row labels.
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row values.
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Congenital fibrosis of the extraocular muscles"
$meta.Range("C2").Value = 512

# data_version must stay textual ("1.12"), not become the number 1.12.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.12"
$meta.Range("D2").Style = "Normal"

$meta.Range("E2").Value = "2021-07-19T10:23:39.378853Z"
$meta.Range("F2").Value = "2021-10-05 14:19:42.636928"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/512/?format=json"

# Leave the "data" sheet selected/active, matching the original workbook.
$data.Select()
